# Auto-generated edit script applying the Unicorn_Profits.xlsx diff
# (sheet names ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR map to sheet1..sheet8)
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 53
$ws.Range("H53").Value = 152.89655
$ws.Range("I53").Value = 108.63636
$ws.Range("J53").Value = 179.94444
$ws.Range("K53").Value = 108.63636
$ws.Range("L53").Value = 179.94444
$ws.Range("M53").Value = 528.36364
$ws.Range("N53").Value = -1453.94444
# Row 98
$ws.Range("H98").Value = 3373.125
$ws.Range("I98").Value = 997.5
$ws.Range("J98").Value = 10500
$ws.Range("K98").Value = 997.5
$ws.Range("L98").Value = 10500
$ws.Range("M98").Value = 500.5
$ws.Range("N98").Value = -13496
# Row 122
$ws.Range("H122").Value = 3373.125
$ws.Range("I122").Value = 997.5
$ws.Range("J122").Value = 10500
$ws.Range("K122").Value = 2992.5
$ws.Range("L122").Value = 31500
$ws.Range("M122").Value = -542.5
$ws.Range("N122").Value = -36400
# Row 135
$ws.Range("H135").Value = 983.6579
$ws.Range("I135").Value = 687.2759
$ws.Range("J135").Value = 1938.6666
$ws.Range("K135").Value = 6185.483099999999
$ws.Range("L135").Value = 17447.9994
$ws.Range("M135").Value = -3650.483099999999
$ws.Range("N135").Value = -22517.9994
# Row 137
$ws.Range("H137").Value = 2369.4883
$ws.Range("I137").Value = 2302.4722
$ws.Range("J137").Value = 2714.1428
$ws.Range("K137").Value = 6907.4166
$ws.Range("L137").Value = 8142.428400000001
$ws.Range("M137").Value = -4357.4166
$ws.Range("N137").Value = -13242.4284
# Row 138
$ws.Range("H138").Value = 2900.5625
$ws.Range("I138").Value = 1967.8667
$ws.Range("J138").Value = 3723.5293
$ws.Range("K138").Value = 5903.6001
$ws.Range("L138").Value = 11170.5879
$ws.Range("M138").Value = -763.6000999999997
$ws.Range("N138").Value = -21450.5879

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 154294.69
$ws.Range("I74").Value = 182675.97
$ws.Range("J74").Value = 62472.94
$ws.Range("K74").Value = 182675.97
$ws.Range("L74").Value = 62472.94
$ws.Range("M74").Value = -181801.97
$ws.Range("N74").Value = -64220.94
# Row 77
$ws.Range("H77").Value = 154294.69
$ws.Range("I77").Value = 182675.97
$ws.Range("J77").Value = 62472.94
$ws.Range("K77").Value = 913379.85
$ws.Range("L77").Value = 312364.7
$ws.Range("M77").Value = -909011.85
$ws.Range("N77").Value = -321100.7
# Row 102
$ws.Range("H102").Value = 2200.5454
$ws.Range("I102").Value = 1336.875
$ws.Range("K102").Value = 1336.875
$ws.Range("M102").Value = 285.125

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 7
$ws.Range("H7").Value = 1630
$ws.Range("I7").Value = 1520
$ws.Range("J7").Value = 1850
$ws.Range("K7").Value = 1520
$ws.Range("L7").Value = 1850
$ws.Range("M7").Value = -1407
$ws.Range("N7").Value = -2076
# Row 141
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").Value = $null

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 5
$ws.Range("H5").Value = 255.2
$ws.Range("I5").Value = 155.5
$ws.Range("K5").Value = 155.5
$ws.Range("M5").Value = -43.5
# Row 25
$ws.Range("H25").Value = 244
$ws.Range("I25").Value = 244
$ws.Range("K25").Value = 244
$ws.Range("M25").Value = -70
# Row 31
$ws.Range("H31").Value = 2650.9312
$ws.Range("I31").Value = 2142.7778
$ws.Range("J31").Value = 3482.4546
$ws.Range("K31").Value = 2142.7778
$ws.Range("L31").Value = 3482.4546
$ws.Range("M31").Value = -1847.7778
$ws.Range("N31").Value = -4072.4546
# Row 34
$ws.Range("H34").Value = 2650.9312
$ws.Range("I34").Value = 2142.7778
$ws.Range("J34").Value = 3482.4546
$ws.Range("K34").Value = 2142.7778
$ws.Range("L34").Value = 3482.4546
$ws.Range("M34").Value = -1940.7778
$ws.Range("N34").Value = -3886.4546
# Row 141
$ws.Range("H141").Value = 29608
$ws.Range("J141").Value = 30469.6
$ws.Range("L141").Value = 30469.6
$ws.Range("N141").Value = -40829.6

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 21
$ws.Range("H21").Value = 750
$ws.Range("J21").Value = 1000
$ws.Range("L21").Value = 3000
$ws.Range("N21").Value = -3346
# Row 34
$ws.Range("H34").Value = 453.625
$ws.Range("J34").Value = 537.8333
$ws.Range("L34").Value = 1613.4999
$ws.Range("N34").Value = -1781.4999
# Row 39
$ws.Range("H39").Value = 5360
$ws.Range("J39").Value = 5600
$ws.Range("L39").Value = 16800
$ws.Range("N39").Value = -17388
# Row 49
$ws.Range("H49").Value = 3000
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").Value = $null
# Row 55
$ws.Range("H55").Value = 2488.889
$ws.Range("J55").Value = 3114.2856
$ws.Range("L55").Value = 9342.856800000001
$ws.Range("N55").Value = -9696.856800000001
# Row 64
$ws.Range("H64").Value = 2954
$ws.Range("I64").Value = 2216
$ws.Range("J64").Value = 3200
$ws.Range("K64").Value = 6648
$ws.Range("L64").Value = 9600
$ws.Range("M64").Value = -6378
$ws.Range("N64").Value = -10140
# Row 67
$ws.Range("H67").Value = 2954
$ws.Range("I67").Value = 2216
$ws.Range("J67").Value = 3200
$ws.Range("K67").Value = 6648
$ws.Range("L67").Value = 9600
$ws.Range("M67").Value = -5712
$ws.Range("N67").Value = -11472
# Row 70
$ws.Range("H70").Value = 5071.2
$ws.Range("I70").Value = 1012
$ws.Range("J70").Value = 5522.222
$ws.Range("K70").Value = 3036
$ws.Range("L70").Value = 16566.666
$ws.Range("M70").Value = -2721
$ws.Range("N70").Value = -17196.666
# Row 73
$ws.Range("H73").Value = 5071.2
$ws.Range("I73").Value = 1012
$ws.Range("J73").Value = 5522.222
$ws.Range("K73").Value = 3036
$ws.Range("L73").Value = 16566.666
$ws.Range("M73").Value = -1944
$ws.Range("N73").Value = -18750.666
# Row 76
$ws.Range("H76").Value = 1850
$ws.Range("I76").Value = 500
$ws.Range("J76").Value = 3200
$ws.Range("K76").Value = 1500
$ws.Range("L76").Value = 9600
$ws.Range("M76").Value = -1117
$ws.Range("N76").Value = -10366
# Row 79
$ws.Range("H79").Value = 1850
$ws.Range("I79").Value = 500
$ws.Range("J79").Value = 3200
$ws.Range("K79").Value = 1500
$ws.Range("L79").Value = 9600
$ws.Range("M79").Value = -174
$ws.Range("N79").Value = -12252
# Row 88
$ws.Range("H88").Value = 5081.6665
$ws.Range("J88").Value = 5081.6665
$ws.Range("L88").Value = 15244.9995
$ws.Range("N88").Value = -16100.9995
# Row 91
$ws.Range("H91").Value = 5081.6665
$ws.Range("J91").Value = 5081.6665
$ws.Range("L91").Value = 15244.9995
$ws.Range("N91").Value = -18208.9995
# Row 92
$ws.Range("H92").Value = 125001064
$ws.Range("I92").Value = 166667120
$ws.Range("J92").Value = 2900
$ws.Range("K92").Value = 500001360
$ws.Range("L92").Value = 8700
$ws.Range("M92").Value = -500000112
$ws.Range("N92").Value = -11196
# Row 94
$ws.Range("H94").Value = 7170.6665
# Row 100
$ws.Range("H100").Value = 6657.609
$ws.Range("J100").Value = 6909.524
$ws.Range("L100").Value = 20728.572
$ws.Range("N100").Value = -22350.572
# Row 103
$ws.Range("H103").Value = 412.14285
$ws.Range("I103").Value = 412.14285
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 1236.42855
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -357.4285500000001
$ws.Range("N103").Value = $null
# Row 106
$ws.Range("H106").Value = 5000
$ws.Range("J106").Value = 5000
$ws.Range("L106").Value = 15000
$ws.Range("N106").Value = -16892

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 4207.5713
$ws.Range("I132").Value = 4480
$ws.Range("J132").Value = 3884.0625
$ws.Range("K132").Value = 13440
$ws.Range("L132").Value = 11652.1875
$ws.Range("M132").Value = -10910
$ws.Range("N132").Value = -16712.1875

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 6058.9614
$ws.Range("I132").Value = 2103.1614
$ws.Range("J132").Value = 11898.477
$ws.Range("K132").Value = 6309.4842
$ws.Range("L132").Value = 35695.431
$ws.Range("M132").Value = -3779.4842
$ws.Range("N132").Value = -40755.431
# Row 136
$ws.Range("H136").Value = 6143.107
$ws.Range("I136").Value = 4277.154
$ws.Range("J136").Value = 7760.2666
$ws.Range("K136").Value = 12831.462
$ws.Range("L136").Value = 23280.7998
$ws.Range("M136").Value = -10281.462
$ws.Range("N136").Value = -28380.7998
